$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$layout = $master.CustomLayouts.Item(2)

# 1. Remove the old "Cloud Foundry" logo picture (Picture 4) from the layout.
$oldLogo = $layout.Shapes.Item("Picture 4")
$oldLogo.Delete()

# 2. Resize + relabel the "SAP HANA Cloud" caption textbox -> "Easy Franchise DB".
$caption = $layout.Shapes.Item("TextBox 29")
$caption.Width = 947695 / 12700
$caption.TextFrame.TextRange.Text = "Easy Franchise DB"
$caption.TextFrame.TextRange.Font.Size = 8
$caption.TextFrame.TextRange.Font.Bold = $true
$caption.TextFrame.TextRange.Font.Color.RGB = 8421504

# 3. Add a new "SAP HANA Cloud" caption textbox (multicloud label) to the layout.
$newBox = $layout.Shapes.AddTextbox(1, 6443895 / 12700, 1067223 / 12700, 2647483 / 12700, 246221 / 12700)
$newBox.Name = "TextBox 57"
$newBox.TextFrame.TextRange.Text = "SAP HANA Cloud"
$newBox.TextFrame.TextRange.Font.Size = 10
$newBox.TextFrame.TextRange.Font.Bold = $true
$newBox.TextFrame.TextRange.Font.Color.RGB = 13998939
